$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# League-bases refresh (28-06-2024 19:47): the match rows for this round got
# re-fetched/re-ordered. Rows 267-271 keep their row id / Div / Date (cols
# A, C, D) but columns B (match id) and E:AD (teams, scores, odds) are
# rewritten with the refreshed data set, which amounts to the following
# rotation of the previous per-row payloads:
#   row 267 <- old row 268
#   row 268 <- old row 267
#   row 269 <- old row 270
#   row 270 <- old row 271
#   row 271 <- old row 269

# Positional params (this host only binds positionally, not by -Name):
# Row, B(id), E(HomeTeam), F(AwayTeam), G(FTHG), H(FTAG), I(HTHG), J(HTAG),
# K(FTR), L(oddH_op), M(oddD_op), N(oddA_op), O(oddH), P(oddD), Q(oddA),
# R(Ah), S(oddAHH), T(oddAHA), U(AhOU), V(oddAHOver), W(oddAHUnder),
# X(PLH), Y(PLD), Z(PLA), AA(PL_Ahh), AB(PL_Aha), AC(PL_AhOver), AD(PL_AhUnder)
function Set-Row {
    param($Row, $B, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U, $V, $W, $X, $Y, $Z, $AA, $AB, $AC, $AD)

    $ws.Cells.Item($Row, 2).Value2 = $B
    $ws.Cells.Item($Row, 5).Value2 = $E
    $ws.Cells.Item($Row, 6).Value2 = $F
    $ws.Cells.Item($Row, 7).Value2 = $G
    $ws.Cells.Item($Row, 8).Value2 = $H
    $ws.Cells.Item($Row, 9).Value2 = $I
    $ws.Cells.Item($Row, 10).Value2 = $J
    $ws.Cells.Item($Row, 11).Value2 = $K
    $ws.Cells.Item($Row, 12).Value2 = $L
    $ws.Cells.Item($Row, 13).Value2 = $M
    $ws.Cells.Item($Row, 14).Value2 = $N
    $ws.Cells.Item($Row, 15).Value2 = $O
    $ws.Cells.Item($Row, 16).Value2 = $P
    $ws.Cells.Item($Row, 17).Value2 = $Q
    $ws.Cells.Item($Row, 18).Value2 = $R
    $ws.Cells.Item($Row, 19).Value2 = $S
    $ws.Cells.Item($Row, 20).Value2 = $T
    $ws.Cells.Item($Row, 21).Value2 = $U
    $ws.Cells.Item($Row, 22).Value2 = $V
    $ws.Cells.Item($Row, 23).Value2 = $W
    $ws.Cells.Item($Row, 24).Value2 = $X
    $ws.Cells.Item($Row, 25).Value2 = $Y
    $ws.Cells.Item($Row, 26).Value2 = $Z
    $ws.Cells.Item($Row, 27).Value2 = $AA
    $ws.Cells.Item($Row, 28).Value2 = $AB
    $ws.Cells.Item($Row, 29).Value2 = $AC
    $ws.Cells.Item($Row, 30).Value2 = $AD
}

Set-Row 267 8203655 "Municipal Perez Zeledon" "Municipal Liberia" `
    0 3 0 1 "A" `
    3.3 3.5 2 2.9 3.4 2.2 0.25 `
    1.825 1.975 2.75 2 1.8 `
    -1 -1 1.2 `
    -1 0.9750000000000001 0.5 -0.5

Set-Row 268 8162891 "Deportivo Saprissa" "Santos de Gupiles" `
    3 1 2 0 "H" `
    1.166 6.5 13 1.125 8 15 -2.25 `
    1.825 1.975 3.5 1.975 1.825 `
    0.125 -1 -1 `
    -0.5 0.4875 0.9750000000000001 -1

Set-Row 269 8162892 "Alajuelense" "AD Guanacasteca" `
    5 0 2 0 "H" `
    1.25 5 10 1.3 4.75 8 -1.5 `
    1.9 1.9 3 1.9 1.9 `
    0.3 -1 -1 `
    0.8999999999999999 -1 0.8999999999999999 -1

Set-Row 270 8162893 "AD Grecia" "AD San Carlos" `
    2 2 0 1 "D" `
    5 4 1.533 4.2 4.2 1.6 1 `
    1.775 2.025 3 1.925 1.875 `
    -1 3.2 -1 `
    0.7749999999999999 -1 0.925 -1

Set-Row 271 8162895 "Sporting San Jose" "Herediano" `
    1 1 1 0 "D" `
    3.6 3.5 1.833 4.5 3.8 1.571 0.75 `
    2.025 1.775 2.75 1.975 1.825 `
    -1 2.8 -1 `
    1.025 -1 -1 0.825
